# Append new sensor-log rows to PIR, Humidity and Temperature sheets
# (auto-logged device readings for 2026-01-28 18:12-18:13).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$rng = $ws.Range("A136:F149")
$arr = New-Object 'object[,]' 14,6
$arr[0,0] = "2026-01-28"; $arr[0,1] = "18:12:44"; $arr[0,2] = "18:00"; $arr[0,3] = "Bathroom"; $arr[0,4] = "No Motion"; $arr[0,5] = "Inactive"
$arr[1,0] = "2026-01-28"; $arr[1,1] = "18:12:47"; $arr[1,2] = "18:00"; $arr[1,3] = "Bathroom"; $arr[1,4] = "No Motion"; $arr[1,5] = "Inactive"
$arr[2,0] = "2026-01-28"; $arr[2,1] = "18:12:49"; $arr[2,2] = "18:00"; $arr[2,3] = "Bathroom"; $arr[2,4] = "No Motion"; $arr[2,5] = "Inactive"
$arr[3,0] = "2026-01-28"; $arr[3,1] = "18:12:52"; $arr[3,2] = "18:00"; $arr[3,3] = "Bathroom"; $arr[3,4] = "No Motion"; $arr[3,5] = "Inactive"
$arr[4,0] = "2026-01-28"; $arr[4,1] = "18:12:57"; $arr[4,2] = "18:00"; $arr[4,3] = "Bathroom"; $arr[4,4] = "No Motion"; $arr[4,5] = "Inactive"
$arr[5,0] = "2026-01-28"; $arr[5,1] = "18:13:01"; $arr[5,2] = "18:00"; $arr[5,3] = "Bathroom"; $arr[5,4] = "No Motion"; $arr[5,5] = "Inactive"
$arr[6,0] = "2026-01-28"; $arr[6,1] = "18:13:05"; $arr[6,2] = "18:00"; $arr[6,3] = "Bathroom"; $arr[6,4] = "No Motion"; $arr[6,5] = "Inactive"
$arr[7,0] = "2026-01-28"; $arr[7,1] = "18:13:10"; $arr[7,2] = "18:00"; $arr[7,3] = "Bathroom"; $arr[7,4] = "No Motion"; $arr[7,5] = "Inactive"
$arr[8,0] = "2026-01-28"; $arr[8,1] = "18:13:17"; $arr[8,2] = "18:00"; $arr[8,3] = "Bathroom"; $arr[8,4] = "No Motion"; $arr[8,5] = "Inactive"
$arr[9,0] = "2026-01-28"; $arr[9,1] = "18:13:21"; $arr[9,2] = "18:00"; $arr[9,3] = "Bathroom"; $arr[9,4] = "No Motion"; $arr[9,5] = "Inactive"
$arr[10,0] = "2026-01-28"; $arr[10,1] = "18:13:25"; $arr[10,2] = "18:00"; $arr[10,3] = "Bathroom"; $arr[10,4] = "No Motion"; $arr[10,5] = "Inactive"
$arr[11,0] = "2026-01-28"; $arr[11,1] = "18:13:30"; $arr[11,2] = "18:00"; $arr[11,3] = "Bathroom"; $arr[11,4] = "No Motion"; $arr[11,5] = "Inactive"
$arr[12,0] = "2026-01-28"; $arr[12,1] = "18:13:35"; $arr[12,2] = "18:00"; $arr[12,3] = "Bathroom"; $arr[12,4] = "No Motion"; $arr[12,5] = "Inactive"
$arr[13,0] = "2026-01-28"; $arr[13,1] = "18:13:41"; $arr[13,2] = "18:00"; $arr[13,3] = "Bathroom"; $arr[13,4] = "No Motion"; $arr[13,5] = "Inactive"
$rng.NumberFormat = "@"
$rng.Value = $arr
$rng.ClearFormats()

$ws = $wb.Worksheets.Item("Humidity")
$rng = $ws.Range("A132:F143")
$arr = New-Object 'object[,]' 12,6
$arr[0,0] = "2026-01-28"; $arr[0,1] = "18:12:45"; $arr[0,2] = "18:00"; $arr[0,3] = "Bathroom"; $arr[0,4] = "88.2%"; $arr[0,5] = "Active"
$arr[1,0] = "2026-01-28"; $arr[1,1] = "18:12:47"; $arr[1,2] = "18:00"; $arr[1,3] = "Bathroom"; $arr[1,4] = "88.2%"; $arr[1,5] = "Active"
$arr[2,0] = "2026-01-28"; $arr[2,1] = "18:12:50"; $arr[2,2] = "18:00"; $arr[2,3] = "Bathroom"; $arr[2,4] = "87.3%"; $arr[2,5] = "Active"
$arr[3,0] = "2026-01-28"; $arr[3,1] = "18:12:53"; $arr[3,2] = "18:00"; $arr[3,3] = "Bathroom"; $arr[3,4] = "88.2%"; $arr[3,5] = "Active"
$arr[4,0] = "2026-01-28"; $arr[4,1] = "18:12:55"; $arr[4,2] = "18:00"; $arr[4,3] = "Bathroom"; $arr[4,4] = "87.3%"; $arr[4,5] = "Active"
$arr[5,0] = "2026-01-28"; $arr[5,1] = "18:12:59"; $arr[5,2] = "18:00"; $arr[5,3] = "Bathroom"; $arr[5,4] = "88.2%"; $arr[5,5] = "Active"
$arr[6,0] = "2026-01-28"; $arr[6,1] = "18:13:03"; $arr[6,2] = "18:00"; $arr[6,3] = "Bathroom"; $arr[6,4] = "88.3%"; $arr[6,5] = "Active"
$arr[7,0] = "2026-01-28"; $arr[7,1] = "18:13:07"; $arr[7,2] = "18:00"; $arr[7,3] = "Bathroom"; $arr[7,4] = "87.3%"; $arr[7,5] = "Active"
$arr[8,0] = "2026-01-28"; $arr[8,1] = "18:13:15"; $arr[8,2] = "18:00"; $arr[8,3] = "Bathroom"; $arr[8,4] = "87.3%"; $arr[8,5] = "Active"
$arr[9,0] = "2026-01-28"; $arr[9,1] = "18:13:19"; $arr[9,2] = "18:00"; $arr[9,3] = "Bathroom"; $arr[9,4] = "88.2%"; $arr[9,5] = "Active"
$arr[10,0] = "2026-01-28"; $arr[10,1] = "18:13:31"; $arr[10,2] = "18:00"; $arr[10,3] = "Bathroom"; $arr[10,4] = "88.2%"; $arr[10,5] = "Active"
$arr[11,0] = "2026-01-28"; $arr[11,1] = "18:13:39"; $arr[11,2] = "18:00"; $arr[11,3] = "Bathroom"; $arr[11,4] = "87.3%"; $arr[11,5] = "Active"
$rng.NumberFormat = "@"
$rng.Value = $arr
$rng.ClearFormats()

$ws = $wb.Worksheets.Item("Temperature")
$rng = $ws.Range("A131:F143")
$arr = New-Object 'object[,]' 13,6
$arr[0,0] = "2026-01-28"; $arr[0,1] = "18:12:43"; $arr[0,2] = "18:00"; $arr[0,3] = "Bathroom"; $arr[0,4] = "22.9C"; $arr[0,5] = "Active"
$arr[1,0] = "2026-01-28"; $arr[1,1] = "18:12:46"; $arr[1,2] = "18:00"; $arr[1,3] = "Bathroom"; $arr[1,4] = "22.9C"; $arr[1,5] = "Active"
$arr[2,0] = "2026-01-28"; $arr[2,1] = "18:12:48"; $arr[2,2] = "18:00"; $arr[2,3] = "Bathroom"; $arr[2,4] = "22.9C"; $arr[2,5] = "Active"
$arr[3,0] = "2026-01-28"; $arr[3,1] = "18:12:51"; $arr[3,2] = "18:00"; $arr[3,3] = "Bathroom"; $arr[3,4] = "22.9C"; $arr[3,5] = "Active"
$arr[4,0] = "2026-01-28"; $arr[4,1] = "18:12:54"; $arr[4,2] = "18:00"; $arr[4,3] = "Bathroom"; $arr[4,4] = "22.9C"; $arr[4,5] = "Active"
$arr[5,0] = "2026-01-28"; $arr[5,1] = "18:12:56"; $arr[5,2] = "18:00"; $arr[5,3] = "Bathroom"; $arr[5,4] = "22.9C"; $arr[5,5] = "Active"
$arr[6,0] = "2026-01-28"; $arr[6,1] = "18:13:00"; $arr[6,2] = "18:00"; $arr[6,3] = "Bathroom"; $arr[6,4] = "22.9C"; $arr[6,5] = "Active"
$arr[7,0] = "2026-01-28"; $arr[7,1] = "18:13:04"; $arr[7,2] = "18:00"; $arr[7,3] = "Bathroom"; $arr[7,4] = "23.0C"; $arr[7,5] = "Active"
$arr[8,0] = "2026-01-28"; $arr[8,1] = "18:13:08"; $arr[8,2] = "18:00"; $arr[8,3] = "Bathroom"; $arr[8,4] = "22.9C"; $arr[8,5] = "Active"
$arr[9,0] = "2026-01-28"; $arr[9,1] = "18:13:16"; $arr[9,2] = "18:00"; $arr[9,3] = "Bathroom"; $arr[9,4] = "22.9C"; $arr[9,5] = "Active"
$arr[10,0] = "2026-01-28"; $arr[10,1] = "18:13:20"; $arr[10,2] = "18:00"; $arr[10,3] = "Bathroom"; $arr[10,4] = "22.9C"; $arr[10,5] = "Active"
$arr[11,0] = "2026-01-28"; $arr[11,1] = "18:13:32"; $arr[11,2] = "18:00"; $arr[11,3] = "Bathroom"; $arr[11,4] = "22.9C"; $arr[11,5] = "Active"
$arr[12,0] = "2026-01-28"; $arr[12,1] = "18:13:40"; $arr[12,2] = "18:00"; $arr[12,3] = "Bathroom"; $arr[12,4] = "22.9C"; $arr[12,5] = "Active"
$rng.NumberFormat = "@"
$rng.Value = $arr
$rng.ClearFormats()

